# The deck currently has the "Integral" theme applied to the slide master
# (ppt/theme/theme2.xml) while the original "Office Theme" colors only
# survive in the notes-master theme (ppt/theme/theme1.xml). The authored
# change swaps the two themes' contents so the slide master goes back to
# the plain "Office Theme" colour scheme.
#
# Re-apply the "Office Theme" colour scheme (theme1.xml's values) onto the
# presentation's active theme colour scheme via the 12-slot ThemeColorScheme
# (Dark1, Light1, Dark2, Light2, Accent1-6, Hyperlink, FollowedHyperlink) -
# the same order/semantics as PowerPoint's MsoThemeColorSchemeIndex.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

$tcs.Item(1).RGB  = 0x000000   # Dark1    -> 000000
$tcs.Item(2).RGB  = 0xFFFFFF   # Light1   -> FFFFFF
$tcs.Item(3).RGB  = 0x6A5444   # Dark2    -> 44546A
$tcs.Item(4).RGB  = 0xE6E6E7   # Light2   -> E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B   # Accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 0x317DED   # Accent2  -> ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5   # Accent3  -> A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF   # Accent4  -> FFC000
$tcs.Item(9).RGB  = 0xC47244   # Accent5  -> 4472C4
$tcs.Item(10).RGB = 0x47AD70   # Accent6  -> 70AD47
$tcs.Item(11).RGB = 0xC16305   # Hyperlink -> 0563C1
$tcs.Item(12).RGB = 0x724F95   # FollowedHyperlink -> 954F72
